$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("K2").Value = 58.5
$ws.Range("N2").Value = 51.15965480231979

# Row 3 updates
$ws.Range("D3").Value = 91159.81
$ws.Range("E3").Value = 62.9
$ws.Range("F3").Value = 0.85
$ws.Range("K3").Value = 53.3
$ws.Range("N3").Value = 51.15965480231979

# Row 4 updates
$ws.Range("K4").Value = 50.3
$ws.Range("N4").Value = 51.15965480231979

# Row 5 updates
$ws.Range("K5").Value = 48.5
$ws.Range("N5").Value = 51.15965480231979

# Row 6 updates
$ws.Range("K6").Value = 34.7
$ws.Range("N6").Value = 51.15965480231979
